# Update cryptos.xlsx with latest price/volume/hour data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells we touch are treated as text so values round-trip exactly
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Cells.Item(2, 4).Value = "306.72"
$ws.Cells.Item(3, 4).Value = "37.29"
$ws.Cells.Item(4, 4).Value = "5.119"
$ws.Cells.Item(5, 4).Value = "0.07759"
$ws.Cells.Item(6, 4).Value = "4.397"
$ws.Cells.Item(7, 4).Value = "1.903"
$ws.Cells.Item(8, 4).Value = "8.207"
$ws.Cells.Item(9, 4).Value = "3.172"
$ws.Cells.Item(10, 4).Value = "0.9183"
$ws.Cells.Item(11, 4).Value = "0.1248"
$ws.Cells.Item(12, 4).Value = "0.1889"
$ws.Cells.Item(13, 4).Value = "0.08712"
$ws.Cells.Item(14, 4).Value = "0.03411"
$ws.Cells.Item(15, 4).Value = "0.09711"
$ws.Cells.Item(16, 4).Value = "0.001368"
$ws.Cells.Item(17, 4).Value = "0.005907"
$ws.Cells.Item(18, 4).Value = "3.550"
$ws.Cells.Item(20, 4).Value = "0.1284"
$ws.Cells.Item(21, 4).Value = "5.026"
$ws.Cells.Item(22, 4).Value = "0.2498"
$ws.Cells.Item(23, 4).Value = "0.02112"
$ws.Cells.Item(24, 4).Value = "0.04338"
$ws.Cells.Item(25, 4).Value = "0.001220"
$ws.Cells.Item(26, 4).Value = "0.004487"
$ws.Cells.Item(27, 4).Value = "0.0001355"
$ws.Cells.Item(39, 4).Value = "0.02199"
$ws.Cells.Item(40, 4).Value = "0.04919"
$ws.Cells.Item(41, 4).Value = "0.007699"
$ws.Cells.Item(42, 4).Value = "0.009831"
$ws.Cells.Item(43, 4).Value = "0.1335"
$ws.Cells.Item(44, 4).Value = "0.002068"
$ws.Cells.Item(45, 4).Value = "0.008779"
$ws.Cells.Item(46, 4).Value = "0.00006870"
$ws.Cells.Item(47, 4).Value = "0.00000000753"
$ws.Cells.Item(48, 4).Value = "0.003012"
$ws.Cells.Item(49, 4).Value = "0.001305"
$ws.Cells.Item(50, 4).Value = "0.00002108"
$ws.Cells.Item(51, 4).Value = "0.0002008"

# --- Volume(1h) (column E) updates ---
$ws.Cells.Item(2, 5).Value = "-3.29%"
$ws.Cells.Item(3, 5).Value = "-6.49%"
$ws.Cells.Item(4, 5).Value = "-0.62%"
$ws.Cells.Item(5, 5).Value = "-5.73%"
$ws.Cells.Item(6, 5).Value = "1.42%"
$ws.Cells.Item(7, 5).Value = "-8.02%"
$ws.Cells.Item(8, 5).Value = "-1.38%"
$ws.Cells.Item(9, 5).Value = "-5.60%"
$ws.Cells.Item(10, 5).Value = "-2.29%"
$ws.Cells.Item(11, 5).Value = "-9.38%"
$ws.Cells.Item(12, 5).Value = "-4.25%"
$ws.Cells.Item(13, 5).Value = "-4.09%"
$ws.Cells.Item(14, 5).Value = "-3.01%"
$ws.Cells.Item(15, 5).Value = "-1.07%"
$ws.Cells.Item(16, 5).Value = "-0.59%"
$ws.Cells.Item(17, 5).Value = "-4.52%"
$ws.Cells.Item(18, 5).Value = "-4.14%"
$ws.Cells.Item(19, 5).Value = "-3.52%"
$ws.Cells.Item(20, 5).Value = "-1.82%"
$ws.Cells.Item(21, 5).Value = "0.56%"
$ws.Cells.Item(22, 5).Value = "1.90%"
$ws.Cells.Item(23, 5).Value = "5,179.18%"
$ws.Cells.Item(24, 5).Value = "-0.33%"
$ws.Cells.Item(25, 5).Value = "-0.97%"
$ws.Cells.Item(26, 5).Value = "-7.01%"
$ws.Cells.Item(27, 5).Value = "4.21%"
$ws.Cells.Item(39, 5).Value = "-0.36%"
$ws.Cells.Item(40, 5).Value = "-5.72%"
$ws.Cells.Item(41, 5).Value = "-0.14%"
$ws.Cells.Item(42, 5).Value = "0.55%"
$ws.Cells.Item(43, 5).Value = "-5.13%"
$ws.Cells.Item(44, 5).Value = "0.88%"
$ws.Cells.Item(45, 5).Value = "-9.08%"
$ws.Cells.Item(46, 5).Value = "3.75%"
$ws.Cells.Item(47, 5).Value = "0.36%"
$ws.Cells.Item(48, 5).Value = "2.39%"
$ws.Cells.Item(49, 5).Value = "-22.77%"
$ws.Cells.Item(50, 5).Value = "0.36%"
$ws.Cells.Item(51, 5).Value = "0.36%"

# --- Hora (column G) updates: all rows 2-51 move from 19 to 20 ---
$ws.Range("G2:G51").Value = "20"

# Restore default (General) styling so the cells keep their original look
$ws.Range("D2:E51").Style = "Normal"
$ws.Range("G2:G51").Style = "Normal"

